# TC09_Canine_StudyUBC02-Breed_Diagnosis_PrimDiseaseSite.xlsx
# "updated ubc2 10 scripts, renamed test suites with w, commiting stashed changes"
#
# Content change: the Cypher query stored in B2 ("CasesTab" query) dropped its
# trailing `co.cohort_description` AS `Cohort` output column (and the now
# -trailing comma on the previous line).
#
# The rest of the diff is view-state noise left behind by the author's Excel
# session (row heights / column widths reflowing after a different monitor
# /DPI, the active selection moving to B2, etc.) which we reproduce as closely
# as the object model allows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- 1. Update the CasesTab query text in B2: remove the trailing Cohort column ---
$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $newQuery

# --- 2. Row heights reflowed (smaller/auto-fit row 2/3/4 wrap heights) ---
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 275.5

# --- 3. Column widths shrank slightly too ---
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 38.83333333333333

# --- 4. Selection / scroll position moved to B2 (with row 1 scrolled out of view) ---
[void]$ws.Range("B2").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
